$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 12: quantity of "sensor ultrasonico" changes from 4 to 1 (total recalculates automatically)
$ws.Range("C12").Value = 1

# Row 15: fill in the previously-empty row with "sensor seguimiento" data
$ws.Range("A15").Value = "sensor seguimiento"
$ws.Range("B15").Value = 1809
$ws.Range("C15").Value = 3

# Update the active selection to A12, matching the saved view state
$ws.Range("A12").Select()
